# Generate Report for Handback
# Updates the localization-status workbook to reflect a failed handback
# transform for the 6b164e09-... file: status text changes on the Overview
# sheet, and per-locale "Error Detail" cells are populated with the
# handback/handoff filename mismatch message. Also widens the Error Detail
# column so the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 6b164e09-... file; its zh-cn/de-de status
# columns (E3/F3) move from "Ready for handoff" to "Handback transform failed".
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# The same shared string backs the per-locale "Status" column (C3) for this
# file on both the zh-cn and de-de sheets, so it flips to the same text too.
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# zh-cn sheet: row 3 (6b164e09-...) Error Detail column (P) gets the
# handback/handoff filename mismatch message.
$zhcn.Range("P3").Value = "Handback file name: 20fof0cl.jpx is different with handoff file name: 6b164e09-6556-4edd-9cd1-a57bc4091425.67560c2c344a67486c5376a6e02582f03c2b3214.zh-cn."

# de-de sheet: row 3 (6b164e09-...) Error Detail column (P) gets the
# handback/handoff filename mismatch message.
$dede.Range("P3").Value = "Handback file name: 20fof0cl.jpx is different with handoff file name: 6b164e09-6556-4edd-9cd1-a57bc4091425.67560c2c344a67486c5376a6e02582f03c2b3214.de-de."

# Widen the Error Detail column (P, the 16th column) on both locale sheets
# so the new, longer message is visible. ColumnWidth is in character units;
# the engine adds the standard ~0.8333 padding when it stores the OOXML
# <col width>, so subtract that to land exactly on width="40".
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666666
$dede.Columns.Item(16).ColumnWidth = 39.16666666666666
